# Add a new header row ("username" / "password") above the existing
# credential row, pushing the existing VIJNARA / Daimler@123 row (and its
# mailto hyperlink) down from row 1 to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing row down to make room for the new header row.
$ws.Rows.Item(1).Insert()

# Write the new header row.
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# The row insert shifts cell contents/styles but leaves the worksheet's
# hyperlink anchored to its original address (B1), so re-point it at the
# credential cell's new home (B2) explicitly: drop the stale hyperlink and
# recreate it there, then restore the "Hyperlink" cell style that Add()
# disturbs.
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Daimler@123") | Out-Null
$ws.Range("B2").Style = "Hyperlink"

# Match the author's final selection (on the now-empty username cell).
$ws.Range("A2").Select() | Out-Null
